# Rewrite the WR_aggregate data: re-sort rows alphabetically by player,
# refresh stat values, and drop the stale trailing 'Mike Evans' rows
# that used to sit at the end before the resort (net -3 rows: 58 -> 55).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Allen Lazard', 'Group1', 13.36666666666667, 9.200000000000001, 63.06666666666666),
    @('Allen Lazard', 'Group2', 13.63333333333333, 7.666666666666667, 51.36666666666667),
    @('Allen Lazard', 'Difference', 0.2666666666666693, -1.533333333333334, -11.7),
    @('Amari Cooper', 'Group1', 13.3, 8.966666666666667, 59.96666666666666),
    @('Amari Cooper', 'Group2', 14.95555555555556, 8.466666666666667, 49.95555555555555),
    @('Amari Cooper', 'Difference', 1.655555555555557, -0.5, -10.01111111111111),
    @('Brandin Cooks', 'Group1', 13.2, 8.499999999999998, 55),
    @('Brandin Cooks', 'Group2', 11.5, 6.8, 50.33333333333334),
    @('Brandin Cooks', 'Difference', -1.700000000000001, -1.699999999999998, -4.666666666666664),
    @('Cedrick Wilson Jr.', '2019-2021', 11.23333333333333, 7.5, 54.76666666666667),
    @('Cedrick Wilson Jr.', '2022-2024', 11.8, 7.733333333333333, 52.2),
    @('Cedrick Wilson Jr.', 'Difference', 0.5666666666666664, 0.2333333333333334, -2.566666666666677),
    @('Darius Slayton', 'Group1', 14.46666666666667, 7.466666666666668, 42.93333333333334),
    @('Darius Slayton', 'Group2', 15.26666666666667, 9.333333333333334, 52.6),
    @('Darius Slayton', 'Difference', 0.7999999999999989, 1.866666666666666, 9.666666666666664),
    @('DeVante Adams', 'Group1', 12.16666666666667, 8.766666666666667, 60.79999999999999),
    @('DeVante Adams', 'Group2', 12.85555555555555, 7.48888888888889, 47.6),
    @('DeVante Adams', 'Difference', 0.6888888888888882, -1.277777777777778, -13.19999999999999),
    @('Diontae Johnson', 'Group1', 10.96666666666667, 6.900000000000001, 47.79999999999999),
    @('Diontae Johnson', 'Group2', 11.075, 6.066666666666666, 44.13333333333333),
    @('Diontae Johnson', 'Difference', 0.1083333333333307, -0.8333333333333348, -3.666666666666657),
    @('joshreynolds', 'Group1', 13.44444444444444, 7.677777777777777, 51.73333333333333),
    @('joshreynolds', 'Group2', 13.84444444444444, 8.077777777777778, 50.62222222222223),
    @('joshreynolds', 'Difference', 0.4000000000000004, 0.4000000000000004, -1.1111111111111),
    @('kalifraymond', 'Group1', 17.23333333333333, 11.6, 56.9),
    @('kalifraymond', 'Group2', 13.23333333333333, 10.16666666666667, 60.43333333333334),
    @('kalifraymond', 'Difference', -4, -1.433333333333334, 3.533333333333339),
    @('Keenan Allen', 'Group1', 10.7, 7.3, 55.63333333333333),
    @('Keenan Allen', 'Group2', 11.16666666666667, 7.600000000000001, 54),
    @('Keenan Allen', 'Difference', 0.4666666666666686, 0.3000000000000016, -1.633333333333333),
    @('Kendrick Bourne', 'Group1', 13.33333333333333, 9.5, 60.96666666666667),
    @('Kendrick Bourne', 'Group2', 11.43333333333333, 8.133333333333333, 52.86666666666667),
    @('Kendrick Bourne', 'Difference', -1.900000000000002, -1.366666666666667, -8.099999999999994),
    @('Marquez Valdes-Scantling', 'Group1', 18.26666666666667, 8.966666666666667, 39.83333333333334),
    @('Marquez Valdes-Scantling', 'Group2', 16.82222222222222, 7.911111111111111, 40.48888888888889),
    @('Marquez Valdes-Scantling', 'Difference', -1.444444444444443, -1.055555555555555, 0.6555555555555515),
    @('Mike Evans', 'Group1', 15.23333333333333, 9.366666666666667, 57.8),
    @('Mike Evans', 'Group2', 14.7, 9.066666666666668, 56.9),
    @('Mike Evans', 'Difference', -0.5333333333333332, -0.2999999999999989, -0.9000000000000057),
    @('Noah Brown', 'Group1', 11.1, 6.866666666666667, 52.06666666666666),
    @('Noah Brown', 'Group2', 14.33333333333333, 8.633333333333333, 52.8),
    @('Noah Brown', 'Difference', 3.233333333333334, 1.766666666666666, 0.7333333333333414),
    @('Ray McCloud', 'Group1', 5.5, 3.85, 37.84999999999999),
    @('Ray McCloud', 'Group2', 13.26666666666667, 8.866666666666667, 56.8),
    @('Ray McCloud', 'Difference', 7.766666666666666, 5.016666666666667, 18.95000000000001),
    @('Tyler Lockett', 'Group1', 13.16666666666667, 9.533333333333333, 61.63333333333335),
    @('Tyler Lockett', 'Group2', 11.93333333333333, 8.066666666666668, 58.56666666666666),
    @('Tyler Lockett', 'Difference', -1.233333333333334, -1.466666666666665, -3.066666666666684),
    @('Tyreek Hill', 'Group1', 13.56666666666667, 9, 59.26666666666666),
    @('Tyreek Hill', 'Group2', 13.76666666666667, 9.466666666666667, 55.46666666666667),
    @('Tyreek Hill', 'Difference', 0.1999999999999975, 0.4666666666666668, -3.79999999999999),
    @('Zach Pascal', 'Group1', 13.06666666666667, 7.633333333333333, 47.53333333333334),
    @('Zach Pascal', 'Group2', 7.4, 4.600000000000001, 41.6),
    @('Zach Pascal', 'Difference', -5.666666666666668, -3.033333333333332, -5.933333333333337),
)

$greenColor = 13434828   # fill used for odd player-groups (style index 2, CCFFCC)
$yellowColor = 12451839  # fill used for even player-groups (style index 3, FFFFBD)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]

    $groupIndex = [Math]::Floor($i / 3)
    if ($groupIndex % 2 -eq 0) {
        $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5)).Interior.Color = $greenColor
    } else {
        $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5)).Interior.Color = $yellowColor
    }
}

# Rows 56:58 held the old trailing 'Mike Evans' block; now redundant
# since Mike Evans was written into rows 38:40 above. Remove them so
# the sheet ends at row 55.
$ws.Range("A56:E58").EntireRow.Delete() | Out-Null
